$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.5679999999999999
$ws.Range("D3").Value = 0.5679999999999999
$ws.Range("D4").Value = 0.264
$ws.Range("D5").Value = 0.538
$ws.Range("D6").Value = 0.497
$ws.Range("D7").Value = 0.314
$ws.Range("D8").Value = 0.32
$ws.Range("D9").Value = 0.271
$ws.Range("D10").Value = 0.497
$ws.Range("D11").Value = 0.311
$ws.Range("D12").Value = 0.32
$ws.Range("D13").Value = 0.271
$ws.Range("D14").Value = 0.497
$ws.Range("D15").Value = 0.314
$ws.Range("D16").Value = 0.32
$ws.Range("D17").Value = 0.271

$ws.Range("C18").Value = 0.67
$ws.Range("D18").Value = 0.616
$ws.Range("D19").Value = 0.629
$ws.Range("D20").Value = 0.617
$ws.Range("D21").Value = 0.61
$ws.Range("D22").Value = 0.626
$ws.Range("D23").Value = 0.618
$ws.Range("D24").Value = 0.615
$ws.Range("D25").Value = 0.613
